# update code tinh luong cho Quyen
# Updates last_edited_time timestamps (column D) for rows 3, 4, 5, 7, 13
# and recalculated "Chi tieu" / "Luy ke" numbers (W5, AA5) for row 5 (Thang 7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# last_edited_time updates -> 2024-07-19T12:51:00.000Z
$ws.Range("D3").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("D4").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("D5").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("D7").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("D13").Value = "2024-07-19T12:51:00.000Z"

# properties.Chi tieu.number (row 5, "Thang 7")
$ws.Range("W5").Value = 17134000

# properties.Luy ke.formula.number (row 5, "Thang 7")
$ws.Range("AA5").Value = 20816000
